$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$delta = [char]0x0394

# Row 9 (A, B, C, D)
$ws.Range("A9").Value = "gdp_q_AA16_obs"
$ws.Range("B9").Value = "Real GDP minus NETEXPORT, net growth, aggregate"
$ws.Range("C9").Value = "Real GDP growth"
$ws.Range("D9").Value = $delta + "LN(GDPCTPI-(NETEXP/GDPTCPI))*100"

# Row 10 (A, B, C) - D10 filled in later
$ws.Range("A10").Value = "i_q_AA16_obs"
$ws.Range("B10").Value = "Norminal investment, net growth, aggrrgate"
$ws.Range("C10").Value = "Nominal investment growth"

# Row 11 (A, C, D) - B11 filled in later
$ws.Range("A11").Value = "c_q_AA16_obs"
$ws.Range("C11").Value = "Consumption growth"
$ws.Range("D11").Value = $delta + "LN(PCESVC96+PCNDGC96)*100"

# Now go back and fill the remaining cells
$ws.Range("D10").Value = $delta + "LN(FPI+PCDGCC96)*100"
$ws.Range("B11").Value = "Consumption, net growth, aggregate"

$ws.Range("B13").Select()
